# Update cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.264.52"
$ws.Range("E2").Value = '  -3.46%  '
$ws.Range("D3").Value = "'2.219.47"
$ws.Range("E3").Value = '  -6.53%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = "'296.95"
$ws.Range("E5").Value = '  -4.40%  '
$ws.Range("D6").Value = "'82.93"
$ws.Range("E6").Value = '  -4.65%  '
$ws.Range("E7").Value = '  -3.26%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  -5.06%  '
$ws.Range("D10").Value = "'0.0775"
$ws.Range("E10").Value = '  -7.87%  '
$ws.Range("D11").Value = "'29.17"
$ws.Range("E11").Value = '  -4.54%  '
$ws.Range("D12").Value = "'47.16"
$ws.Range("E12").Value = '  -10.40%  '
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("D14").Value = "'2.564.99"
$ws.Range("E14").Value = '  -6.08%  '
$ws.Range("D15").Value = "'6.21"
$ws.Range("E15").Value = '  -5.11%  '
$ws.Range("D16").Value = "'14.09"
$ws.Range("E16").Value = '  -6.35%  '
$ws.Range("D17").Value = "'2.229.62"
$ws.Range("E17").Value = '  -6.08%  '
$ws.Range("E18").Value = '  -5.68%  '
$ws.Range("D19").Value = "'39.180.35"
$ws.Range("E19").Value = '  -3.32%  '
$ws.Range("D20").Value = "'0.0₃0871"
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("D21").Value = "'5.71"
$ws.Range("E21").Value = '  -6.77%  '
$ws.Range("D22").Value = "'64.81"
$ws.Range("E22").Value = '  -5.39%  '
$ws.Range("D23").Value = "'10.22"
$ws.Range("E23").Value = '  -4.64%  '
$ws.Range("D24").Value = "'227.02"
$ws.Range("E24").Value = '  -3.33%  '
$ws.Range("E26").Value = '  -7.30%  '
$ws.Range("D27").Value = "'1.76"
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").Value = "'22.60"
$ws.Range("E28").Value = '  -4.42%  '
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").Value = "'32.01"
$ws.Range("E31").Value = '  -5.34%  '
$ws.Range("D32").Value = "'147.73"
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = '  -7.32%  '
$ws.Range("D35").Value = "'0.0697"
$ws.Range("E35").Value = '  -4.36%  '
$ws.Range("E36").Value = '  -4.42%  '
$ws.Range("D37").Value = "'0.109"
$ws.Range("E37").Value = '  -3.75%  '
$ws.Range("D38").Value = "'2.67"
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("D39").Value = "'0.0958"
$ws.Range("E39").Value = '  -3.74%  '
$ws.Range("D40").Value = "'14.84"
$ws.Range("E40").Value = '  -6.38%  '
$ws.Range("E41").Value = '  -4.68%  '
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = '  -3.47%  '
$ws.Range("D43").Value = "'1.913.24"
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("E44").Value = '  -4.09%  '
$ws.Range("D45").Value = "'2.02"
$ws.Range("E45").Value = '  -15.29%  '
$ws.Range("D46").Value = "'8.98"
$ws.Range("E46").Value = '  -4.93%  '
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("D48").Value = "'15.87"
$ws.Range("E48").Value = '  -9.72%  '
$ws.Range("D49").Value = "'2.432.99"
$ws.Range("E49").Value = '  -6.23%  '
$ws.Range("D50").Value = "'70.23"
$ws.Range("E50").Value = '  -2.87%  '
$ws.Range("D51").Value = "'86.91"
$ws.Range("E51").Value = '  -6.76%  '
